# issue #5: stock data from json to db
# The "股票" (stock) sheet gains three new columns: category, source_file, index.

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)   # 股票 (stock) sheet

# 1) Insert a new column at I. This pushes the existing date / legislator_name /
#    legislator_id columns from I,J,K to J,K,L (and carries their formatting along).
$ws4.Columns.Item(9).Insert()
$ws4.Range("I1").Value = "category"

# 2) Append two new trailing columns: source_file, index.
$ws4.Range("M1").Value = "source_file"
$ws4.Range("N1").Value = "index"

# Match the header formatting (bold, centered, bordered) already used by the sheet.
$headerCells = @("M1", "N1")
foreach ($addr in $headerCells) {
    $cell = $ws4.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlHAlignCenter
    $cell.VerticalAlignment = -4160     # xlVAlignTop
    $cell.Borders.LineStyle = 1
}

# 3) Fill in the data rows (2-8): category = "normal", source_file = "tmpf4561",
#    index = same value as column A (the row id).
for ($r = 2; $r -le 8; $r++) {
    $ws4.Range("I$r").Value = "normal"
    $ws4.Range("M$r").Value = "tmpf4561"
    $ws4.Range("N$r").Value = $ws4.Range("A$r").Value2
}
